$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps its original text formatting so purely
# numeric-looking strings (e.g. "0.9997", "1.0000") are not coerced into
# Excel numbers, matching the source workbook where every Price cell is
# stored as text (inline string).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.354.73'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').Value = '1.826.36'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('D4').Value = '0.9997'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '314.83'
$ws.Range('E5').Value = '  +1.56%  '
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.4719'
$ws.Range('E7').Value = '  +6.10%  '
$ws.Range('E8').Value = '  +3.42%  '
$ws.Range('D9').Value = '0.07427'
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('D10').Value = '0.8771'
$ws.Range('E10').Value = '  +2.46%  '
$ws.Range('D11').Value = '20.80'
$ws.Range('E11').Value = '  +0.93%  '
$ws.Range('D12').Value = '1.827.93'
$ws.Range('E12').Value = '  -2.12%  '
$ws.Range('D13').Value = '6.703'
$ws.Range('E13').Value = '  +1.64%  '
$ws.Range('D14').Value = '5.439'
$ws.Range('E14').Value = '  +2.83%  '
$ws.Range('D15').Value = '93.13'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '0.07089'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('D18').Value = '0.000008805'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = '15.03'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').Value = '27.349.01'
$ws.Range('D22').Value = '5.333'
$ws.Range('E22').Value = '  +3.70%  '
$ws.Range('D23').Value = '10.95'
$ws.Range('E23').Value = '  +1.44%  '
$ws.Range('D24').Value = '2.049.80'
$ws.Range('E24').Value = '  -3.11%  '
$ws.Range('D25').Value = '1.939'
$ws.Range('E25').Value = '  -2.63%  '
$ws.Range('D26').Value = '151.13'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('D27').Value = '2.256'
$ws.Range('E27').Value = '  +3.99%  '
$ws.Range('D28').Value = '18.61'
$ws.Range('E28').Value = '  +1.11%  '
$ws.Range('D29').Value = '5.344'
$ws.Range('E29').Value = '  +2.88%  '
$ws.Range('D30').Value = '117.25'
$ws.Range('E30').Value = '  +0.50%  '
$ws.Range('D31').Value = '0.08954'
$ws.Range('E31').Value = '  +1.50%  '
$ws.Range('D32').Value = '0.7945'
$ws.Range('E32').Value = '  +6.52%  '
$ws.Range('D33').Value = '1.194'
$ws.Range('E33').Value = '  +2.01%  '
$ws.Range('D34').Value = '4.551'
$ws.Range('E34').Value = '  +2.20%  '
$ws.Range('D35').Value = '2.940'
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D38').Value = '0.01978'
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('D39').Value = '0.05252'
$ws.Range('E39').Value = '  +1.28%  '
$ws.Range('D40').Value = '7.309'
$ws.Range('E40').Value = '  +4.18%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.5339'
$ws.Range('E41').Value = '  +0.54%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '2.373'
$ws.Range('E42').Value = '  +20.47%  '
$ws.Range('D43').Value = '2.893'
$ws.Range('E43').Value = '  +1.19%  '
$ws.Range('D44').Value = '0.1704'
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('D45').Value = '8.650'
$ws.Range('E45').Value = '  +2.33%  '
$ws.Range('D46').Value = '0.5109'
$ws.Range('E46').Value = '  -0.25%  '
$ws.Range('D47').Value = '10.62'
$ws.Range('E47').Value = '  +1.10%  '
$ws.Range('D48').Value = '105.52'
$ws.Range('E48').Value = '  +0.08%  '
$ws.Range('D49').Value = '1.683'
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('D50').Value = '0.9993'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('D51').Value = '0.06388'
$ws.Range('E51').Value = '  +1.14%  '
